# Update crypto price/volume table with refreshed values from the
# Thu Mar 21 01:56:52 UTC 2024 GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'" + "67.760.52"
$ws.Range("E2").Value = "'" + "  +8.48%  "

# Row 3
$ws.Range("D3").Value = "'" + "3.519.89"
$ws.Range("E3").Value = "'" + "  +10.49%  "

# Row 4
$ws.Range("E4").Value = "'" + "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "'" + "191.00"
$ws.Range("E5").Value = "'" + "  +11.11%  "

# Row 6
$ws.Range("D6").Value = "'" + "552.30"
$ws.Range("E6").Value = "'" + "  +8.21%  "

# Row 7
$ws.Range("D7").Value = "'" + "3.514.95"
$ws.Range("E7").Value = "'" + "  +10.34%  "

# Row 8
$ws.Range("E8").Value = "'" + "  +3.98%  "

# Row 9
$ws.Range("E9").Value = "'" + "  +0.00%  "

# Row 10
$ws.Range("D10").Value = "'" + "0.639"
$ws.Range("E10").Value = "'" + "  +8.28%  "

# Row 11
$ws.Range("D11").Value = "'" + "57.00"
$ws.Range("E11").Value = "'" + "  +9.38%  "

# Row 12
$ws.Range("D12").Value = "'" + "0.151"
$ws.Range("E12").Value = "'" + "  +18.94%  "

# Row 13
$ws.Range("D13").Value = "'" + "0.0000273"
$ws.Range("E13").Value = "'" + "  +10.12%  "

# Row 14
$ws.Range("D14").Value = "'" + "9.46"
$ws.Range("E14").Value = "'" + "  +7.68%  "

# Row 15
$ws.Range("D15").Value = "'" + "4.081.65"
$ws.Range("E15").Value = "'" + "  +10.47%  "

# Row 16
$ws.Range("D16").Value = "'" + "3.518.79"
$ws.Range("E16").Value = "'" + "  +10.92%  "

# Row 17
$ws.Range("D17").Value = "'" + "67.903.61"
$ws.Range("E17").Value = "'" + "  +8.96%  "

# Row 18
$ws.Range("E18").Value = "'" + "  +6.31%  "

# Row 19
$ws.Range("E19").Value = "'" + "  +8.34%  "

# Row 20
$ws.Range("D20").Value = "'" + "11.86"
$ws.Range("E20").Value = "'" + "  +10.26%  "

# Row 21
$ws.Range("E21").Value = "'" + "  +6.22%  "

# Row 22
$ws.Range("D22").Value = "'" + "409.48"
$ws.Range("E22").Value = "'" + "  +13.92%  "

# Row 23
$ws.Range("D23").Value = "'" + "3.94"
$ws.Range("E23").Value = "'" + "  +7.65%  "

# Row 24
$ws.Range("D24").Value = "'" + "84.68"
$ws.Range("E24").Value = "'" + "  +6.92%  "

# Row 25
$ws.Range("E25").Value = "'" + "  +12.62%  "

# Row 26
$ws.Range("D26").Value = "'" + "11.47"
$ws.Range("E26").Value = "'" + "  +4.69%  "

# Row 27
$ws.Range("D27").Value = "'" + "2.99"
$ws.Range("E27").Value = "'" + "  +16.05%  "

# Row 28
$ws.Range("E28").Value = "'" + "  +4.14%  "

# Row 29
$ws.Range("D29").Value = "'" + "12.02"
$ws.Range("E29").Value = "'" + "  +8.89%  "

# Row 30
$ws.Range("D30").Value = "'" + "8.73"
$ws.Range("E30").Value = "'" + "  +8.89%  "

# Row 31
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").Value = "'" + "695.55"
$ws.Range("E31").Value = "'" + "  +8.34%  "

# Row 32
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'" + "30.56"
$ws.Range("E32").Value = "'" + "  +9.63%  "

# Row 33
$ws.Range("D33").Value = "'" + "6.86"
$ws.Range("E33").Value = "'" + "  +10.91%  "

# Row 34
$ws.Range("D34").Value = "'" + "11.81"
$ws.Range("E34").Value = "'" + "  +7.80%  "

# Row 35
$ws.Range("E35").Value = "'" + "  +9.03%  "

# Row 36
$ws.Range("D36").Value = "'" + "60.31"
$ws.Range("E36").Value = "'" + "  +5.75%  "

# Row 37
$ws.Range("D37").Value = "'" + "39.25"
$ws.Range("E37").Value = "'" + "  +9.28%  "

# Row 38
$ws.Range("D38").Value = "'" + "0.0₃0828"
$ws.Range("E38").Value = "'" + "  +21.98%  "

# Row 39
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").Value = "'" + "1.00"
$ws.Range("E39").Value = "'" + "  -0.11%  "

# Row 40
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "'" + "0.398"
$ws.Range("E40").Value = "'" + "  +8.14%  "

# Row 41
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'" + "0.136"
$ws.Range("E41").Value = "'" + "  +12.73%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'" + "3.39"
$ws.Range("E42").Value = "'" + "  +24.88%  "

# Row 43
$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D43").Value = "'" + "3.03"
$ws.Range("E43").Value = "'" + "  +19.52%  "

# Row 44
$ws.Range("D44").Value = "'" + "0.999"
$ws.Range("E44").Value = "'" + "  +0.12%  "

# Row 45
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "'" + "2.71"
$ws.Range("E45").Value = "'" + "  +8.77%  "

# Row 46
$ws.Range("D46").Value = "'" + "3.033.82"
$ws.Range("E46").Value = "'" + "  +7.48%  "

# Row 47
$ws.Range("D47").Value = "'" + "3.40"
$ws.Range("E47").Value = "'" + "  +19.62%  "

# Row 48
$ws.Range("D48").Value = "'" + "0.0422"
$ws.Range("E48").Value = "'" + "  +10.66%  "

# Row 49
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "'" + "2.75"
$ws.Range("E49").Value = "'" + "  +3.18%  "

# Row 50
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "'" + "9.01"
$ws.Range("E50").Value = "'" + "  +22.08%  "

# Row 51
$ws.Range("E51").Value = "'" + "  +7.92%  "

